# NationalTreeSales.xlsx update: "2027 info"
# The columns H, I, J (for every row, header included) get left-rotated:
#   new H = old I
#   new I = old J
#   new J = old H
# Also update the view: scroll so column D is the left-most visible column,
# and select cell J1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NationalTreeSales")

for ($r = 1; $r -le 18; $r++) {
    $oldH = $ws.Cells.Item($r, 8).Value2
    $oldI = $ws.Cells.Item($r, 9).Value2
    $oldJ = $ws.Cells.Item($r, 10).Value2

    $ws.Cells.Item($r, 8).Value2  = $oldI
    $ws.Cells.Item($r, 9).Value2  = $oldJ
    $ws.Cells.Item($r, 10).Value2 = $oldH
}

# Update the view: scroll window so column D is the left-most visible column
# and select J1 as the active cell.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("J1").Select()
